$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(74, 8).Value = 12890.4
$ws.Cells.Item(74, 9).Value = 3613
$ws.Cells.Item(74, 11).Value = 3613
$ws.Cells.Item(74, 13).Value = -2677
$ws.Cells.Item(77, 8).Value = 12890.4
$ws.Cells.Item(77, 9).Value = 3613
$ws.Cells.Item(77, 11).Value = 18065
$ws.Cells.Item(77, 13).Value = -13385
$ws.Cells.Item(80, 8).Value = 300
$ws.Cells.Item(80, 9).Value = 100
$ws.Cells.Item(80, 10).Value = 500
$ws.Cells.Item(80, 11).Value = 300
$ws.Cells.Item(80, 12).Value = 1500
$ws.Cells.Item(80, 13).Value = 698
$ws.Cells.Item(80, 14).Value = -3496
$ws.Cells.Item(83, 8).Value = 300
$ws.Cells.Item(83, 9).Value = 100
$ws.Cells.Item(83, 10).Value = 500
$ws.Cells.Item(83, 11).Value = 900
$ws.Cells.Item(83, 12).Value = 4500
$ws.Cells.Item(83, 13).Value = 4092
$ws.Cells.Item(83, 14).Value = -14484
$ws.Cells.Item(99, 8).Value = 6856.9
$ws.Cells.Item(99, 9).Value = 4849.5
$ws.Cells.Item(99, 11).Value = 14548.5
$ws.Cells.Item(99, 13).Value = -13050.5
$ws.Cells.Item(100, 8).Value = 2520.8
$ws.Cells.Item(100, 9).Value = 2902
$ws.Cells.Item(100, 10).Value = 996
$ws.Cells.Item(100, 11).Value = 2902
$ws.Cells.Item(100, 12).Value = 996
$ws.Cells.Item(100, 13).Value = -2361
$ws.Cells.Item(100, 14).Value = -2078
$ws.Cells.Item(112, 8).Value = 2572
$ws.Cells.Item(112, 10).Value = 1394
$ws.Cells.Item(112, 12).Value = 4182
$ws.Cells.Item(112, 14).Value = -6398
$ws.Cells.Item(129, 8).Value = 2912
$ws.Cells.Item(129, 10).Value = 2873
$ws.Cells.Item(129, 12).Value = 8619
$ws.Cells.Item(129, 14).Value = -18619
$ws.Cells.Item(137, 8).Value = 2336.4167
$ws.Cells.Item(137, 9).Value = 864.5
$ws.Cells.Item(137, 10).Value = 3808.3333
$ws.Cells.Item(137, 11).Value = 2593.5
$ws.Cells.Item(137, 12).Value = 11424.9999
$ws.Cells.Item(137, 13).Value = -43.5
$ws.Cells.Item(137, 14).Value = -16524.9999

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(2, 8).Value = 1089.1666
$ws.Cells.Item(2, 9).Value = 1006.36365
$ws.Cells.Item(2, 11).Value = 1006.36365
$ws.Cells.Item(2, 13).Value = -893.36365
$ws.Cells.Item(44, 8).Value = 14502
$ws.Cells.Item(44, 10).Value = 14502
$ws.Cells.Item(44, 12).Value = 14502
$ws.Cells.Item(44, 14).Value = -15478
$ws.Cells.Item(45, 8).Value = 2993.3076
$ws.Cells.Item(45, 9).Value = 1612.5714
$ws.Cells.Item(45, 11).Value = 1612.5714
$ws.Cells.Item(45, 13).Value = -1235.5714
$ws.Cells.Item(88, 8).Value = 1883.75
$ws.Cells.Item(88, 9).Value = 1506.3334
$ws.Cells.Item(88, 10).Value = 3016
$ws.Cells.Item(88, 11).Value = 1506.3334
$ws.Cells.Item(88, 12).Value = 3016
$ws.Cells.Item(88, 13).Value = -1100.3334
$ws.Cells.Item(88, 14).Value = -3828
$ws.Cells.Item(91, 8).Value = 1883.75
$ws.Cells.Item(91, 9).Value = 1506.3334
$ws.Cells.Item(91, 10).Value = 3016
$ws.Cells.Item(91, 11).Value = 1506.3334
$ws.Cells.Item(91, 12).Value = 3016
$ws.Cells.Item(91, 13).Value = -102.3334
$ws.Cells.Item(91, 14).Value = -5824
$ws.Cells.Item(116, 8).Value = 1089.1666
$ws.Cells.Item(116, 9).Value = 1006.36365
$ws.Cells.Item(116, 11).Value = 1006.36365
$ws.Cells.Item(116, 13).Value = 1287.63635

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(3, 8).Value = 1089.1666
$ws.Cells.Item(3, 9).Value = 1006.36365
$ws.Cells.Item(3, 11).Value = 1006.36365
$ws.Cells.Item(3, 13).Value = -892.36365
$ws.Cells.Item(82, 8).Value = 19787.4
$ws.Cells.Item(82, 9).Value = 5962.875
$ws.Cells.Item(82, 10).Value = 75085.5
$ws.Cells.Item(82, 11).Value = 5962.875
$ws.Cells.Item(82, 12).Value = 75085.5
$ws.Cells.Item(82, 13).Value = -5579.875
$ws.Cells.Item(82, 14).Value = -75851.5
$ws.Cells.Item(85, 8).Value = 19787.4
$ws.Cells.Item(85, 9).Value = 5962.875
$ws.Cells.Item(85, 10).Value = 75085.5
$ws.Cells.Item(85, 11).Value = 5962.875
$ws.Cells.Item(85, 12).Value = 75085.5
$ws.Cells.Item(85, 13).Value = -4636.875
$ws.Cells.Item(85, 14).Value = -77737.5
$ws.Cells.Item(86, 8).Value = 9200
$ws.Cells.Item(86, 9).Value = 0
$ws.Cells.Item(86, 11).Value = 0
$ws.Cells.Item(86, 13).ClearContents()
$ws.Cells.Item(89, 8).Value = 9200
$ws.Cells.Item(89, 9).Value = 0
$ws.Cells.Item(89, 11).Value = 6500
$ws.Cells.Item(89, 13).ClearContents()
$ws.Cells.Item(107, 8).Value = 8700.583000000001
$ws.Cells.Item(107, 9).Value = 8113.375
$ws.Cells.Item(107, 11).Value = 8113.375
$ws.Cells.Item(107, 13).Value = -6193.375

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(58, 8).Value = 3437.111
$ws.Cells.Item(58, 9).Value = 1905.6666
$ws.Cells.Item(58, 11).Value = 1905.6666
$ws.Cells.Item(58, 13).Value = -1702.6666
$ws.Cells.Item(99, 8).Value = 14995
$ws.Cells.Item(99, 10).Value = 14990
$ws.Cells.Item(99, 12).Value = 14990
$ws.Cells.Item(99, 14).Value = -17986
$ws.Cells.Item(105, 8).Value = 560.5
$ws.Cells.Item(105, 9).Value = 510
$ws.Cells.Item(105, 10).Value = 611
$ws.Cells.Item(105, 11).Value = 510
$ws.Cells.Item(105, 12).Value = 611
$ws.Cells.Item(105, 13).Value = 1237
$ws.Cells.Item(105, 14).Value = -4105
$ws.Cells.Item(122, 8).Value = 2997.5
$ws.Cells.Item(122, 9).Value = 2997.5
$ws.Cells.Item(122, 11).Value = 8992.5
$ws.Cells.Item(122, 13).Value = -6542.5
$ws.Cells.Item(126, 8).Value = 14995
$ws.Cells.Item(126, 10).Value = 14990
$ws.Cells.Item(126, 12).Value = 44970
$ws.Cells.Item(126, 14).Value = -49910
$ws.Cells.Item(136, 8).Value = 3437.111
$ws.Cells.Item(136, 9).Value = 1905.6666
$ws.Cells.Item(136, 11).Value = 5716.9998
$ws.Cells.Item(136, 13).Value = -3166.9998

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(11, 8).Value = 871.9375
$ws.Cells.Item(11, 9).Value = 329.41666
$ws.Cells.Item(11, 11).Value = 988.2499799999999
$ws.Cells.Item(11, 13).Value = -848.2499799999999
$ws.Cells.Item(26, 8).Value = 245
$ws.Cells.Item(26, 9).Value = 245
$ws.Cells.Item(26, 11).Value = 735
$ws.Cells.Item(26, 13).Value = -447
$ws.Cells.Item(34, 8).Value = 1505.9375
$ws.Cells.Item(34, 10).Value = 1136.4166
$ws.Cells.Item(34, 12).Value = 3409.2498
$ws.Cells.Item(34, 14).Value = -3577.2498
$ws.Cells.Item(52, 8).Value = 1198.6
$ws.Cells.Item(52, 10).Value = 1198.6
$ws.Cells.Item(52, 12).Value = 3595.8
$ws.Cells.Item(52, 14).Value = -4127.799999999999
$ws.Cells.Item(104, 8).Value = 9028.177
$ws.Cells.Item(104, 10).Value = 9712.786
$ws.Cells.Item(104, 12).Value = 29138.358
$ws.Cells.Item(104, 14).Value = -34380.358
$ws.Cells.Item(109, 8).Value = 2512750
$ws.Cells.Item(109, 9).Value = 6667333.5
$ws.Cells.Item(109, 11).Value = 20002000.5
$ws.Cells.Item(109, 13).Value = -20000960.5
$ws.Cells.Item(132, 8).Value = 2893.4707
$ws.Cells.Item(132, 9).Value = 2676.3845
$ws.Cells.Item(132, 11).Value = 24087.4605
$ws.Cells.Item(132, 13).Value = -21557.4605
$ws.Cells.Item(140, 8).Value = 2232.625
$ws.Cells.Item(140, 9).Value = 2048.1333
$ws.Cells.Item(140, 11).Value = 6144.3999
$ws.Cells.Item(140, 13).Value = -964.3999000000003

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(133, 8).Value = 105666.664
$ws.Cells.Item(133, 10).Value = 105666.664
$ws.Cells.Item(133, 12).Value = 105666.664
$ws.Cells.Item(133, 14).Value = -115786.664

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(122, 8).Value = 2252.6667
$ws.Cells.Item(122, 9).Value = 2252.6667
$ws.Cells.Item(122, 11).Value = 6758.000100000001
$ws.Cells.Item(122, 13).Value = -4308.000100000001

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(54, 8).Value = 88998
$ws.Cells.Item(54, 10).Value = 88998
$ws.Cells.Item(54, 12).Value = 88998
$ws.Cells.Item(54, 14).Value = -90038
$ws.Cells.Item(107, 8).Value = 536.4
$ws.Cells.Item(107, 9).Value = 545.5
$ws.Cells.Item(107, 11).Value = 1636.5
$ws.Cells.Item(107, 13).Value = 283.5
